$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Build the new "Periodo / Pronóstico / Límite inferior / Límite superior"
#    table at rows 59-65 by cloning the style/layout of the existing table
#    at rows 51-57 (same columns C:F, same per-cell styles).
# ---------------------------------------------------------------------------
$srcTable = $ws.Range("C51:F57")
$dstTable = $ws.Range("C59")
$srcTable.Copy($dstTable)

# Row 59 is the header row - give it its own (slightly shorter) custom height.
$ws.Rows.Item(59).RowHeight = 17.5

# Fix up the header text that differs from the source table.
$ws.Range("C59").Value = "Periodo"
$ws.Range("D59").Value = "Pronóstico"
$ws.Range("E59").Value = "Límite inferior"
$ws.Range("F59").Value = "Límite superior"

# ---------------------------------------------------------------------------
# 2) Write the new forecast values for rows 60-65. These values must be
#    stored as *text* (matching the existing table's cells, which are text
#    too) rather than as numbers, so a direct .Value assignment (which would
#    auto-coerce a numeric-looking string into a real number and could drop
#    a trailing zero) cannot be used. Instead we stage each value in a
#    scratch cell as a formula that evaluates to a string, copy *values
#    only* into the destination (this preserves the destination's existing
#    number format / style while still landing a text cell), then clear the
#    scratch cell.
# ---------------------------------------------------------------------------
function Set-TextValue($cellAddress, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy() | Out-Null
    $ws.Range($cellAddress).PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.Clear() | Out-Null
}

Set-TextValue "D60" "10866.35"
Set-TextValue "E60" "10529.82"
Set-TextValue "F60" "11202.87"

Set-TextValue "D61" "10884.31"
Set-TextValue "E61" "10495.27"
Set-TextValue "F61" "11273.35"

Set-TextValue "D62" "10902.10"
Set-TextValue "E62" "10491.03"
Set-TextValue "F62" "11313.16"

Set-TextValue "D63" "10919.82"
Set-TextValue "E63" "10495.52"
Set-TextValue "F63" "11344.12"

Set-TextValue "D64" "10937.50"
Set-TextValue "E64" "10503.14"
Set-TextValue "F64" "11371.87"

Set-TextValue "D65" "10955.17"
Set-TextValue "E65" "10512.07"
Set-TextValue "F65" "11398.28"

# ---------------------------------------------------------------------------
# 3) Column F got narrower (it used to hold long text, now it holds short
#    numbers).
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 13.8333333333333

# ---------------------------------------------------------------------------
# 4) Scroll position / selection: the sheet view now shows row 47 at the top
#    with G56 selected (previously F79 / M93).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G56").Select() | Out-Null
